$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @(99.0, 600.0652608881426, 583017857.9747233, 117354.3013491844, 1206.0652608878343),
        @(5148.0, 600.6952793721468, 443594728.58910286, 88761.41578616571, 1806.6952793713717),
        @(10197.0, 599.80445441261, 13728363014.986067, 4078626.209533501, 2405.8044544113263),
        @(15246.0, 599.5418247920735, 25361699956.151226, 7571742.772391258, 3005.541824790268),
        @(20295.0, 600.1350160793911, 36345058456.54646, 10868401.223173268, 3606.135016077064),
        @(25344.0, 598.7055398472703, 47002255979.12908, 14065855.314072514, 4204.705539845095),
        @(30393.0, 599.5559627779041, 57477999457.25365, 17207904.488108765, 4805.555962779189),
        @(35442.0, 600.3149339394379, 67833971578.14027, 20315779.59313609, 5406.314933944289),
        @(40491.0, 600.7337522650741, 78125627897.65572, 23400049.15706407, 6006.733752273604),
        @(45540.0, 601.0055847436566, 88342696984.0267, 26467442.75678684, 6607.005584755917),
        @(50589.0, 600.0041040861327, 98505097936.01735, 29517460.612074636, 7206.004104102185),
        @(55638.0, 599.6057995355509, 108652442511.21083, 32558333.020670354, 7805.605799555407),
        @(60687.0, 598.8607516553421, 118759174256.89987, 35592880.96811674, 8404.860751676044),
        @(65736.0, 600.7034331867904, 128896361463.1997, 38635090.199479796, 9006.70343320277),
        @(70785.0, 599.942214056184, 138976443399.69513, 41657915.0389202, 9605.942214067545),
        @(75834.0, 599.6454242263169, 149057720658.26907, 44681901.341159455, 10205.64542423293),
        @(80883.0, 600.846133802974, 159150519711.0038, 47708955.10131971, 10806.8461338048),
        @(85932.0, 599.1580724506541, 169184415583.78174, 50719219.26227313, 11405.158072447623),
        @(90981.0, 600.5010047095185, 179273021909.38052, 53744011.055407576, 12006.50100470158),
        @(96030.0, 599.0627416674183, 189301682840.64996, 56754135.770848505, 12605.06274165467),
        @(101079.0, 600.7952934445273, 199373618722.93405, 59775667.64149743, 13206.795293426894),
        @(106128.0, 600.6367534726522, 209404041623.851, 62784085.86919422, 13806.636753450059),
        @(111177.0, 601.0281536276641, 219457574919.8841, 65798595.34928951, 14407.028153600262),
        @(116226.0, 599.8424317561097, 229463371275.18942, 68801761.64976402, 15005.842431723688),
        @(121275.0, 601.0324870683486, 239533227288.2953, 71820772.42951734, 15607.032487030992),
        @(126324.0, 600.3077888599463, 249539728244.98907, 74823766.51168096, 16206.307788817629),
        @(131373.0, 599.9372170605619, 259559636790.81094, 77830651.3089226, 16805.93721702558),
        @(136422.0, 600.2336746771521, 269592611961.22537, 80839420.23262672, 17406.233674654824),
        @(141471.0, 599.7891603848306, 279619438206.7032, 83845693.21045086, 18005.789160375058),
        @(146520.0, 600.235816586224, 289656156811.7445, 86856089.9255712, 18606.235816589167)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $rowVals[$j]
    }
}
